$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("general")

# --- Copy the cell formats that the new rows should inherit -----------------
# Column A on these new rows uses the same shaded/bordered style as C25
# (style index 67 in the original file); column D uses the same style as
# D19 (style index 63 in the original file).
$ws.Range("C25").Copy() | Out-Null
$ws.Range("A26:A28").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("D19").Copy() | Out-Null
$ws.Range("D26:D28").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# --- New GPS offset rows (26-28) --------------------------------------------
# Column A: parameter names
$ws.Range("A26").Value = "r_gps_x"
$ws.Range("A27").Value = "r_gps_y"
$ws.Range("A28").Value = "r_gps_z"

# Column B: values
$ws.Range("B26").Value = 0
$ws.Range("B27").Value = -0.25
$ws.Range("B28").Value = 0

# Column D: descriptions
$ws.Range("D26").Value = "x component of gps position wrt to body frame"
$ws.Range("D27").Value = "y component of gps position wrt to body frame"
$ws.Range("D28").Value = "z component of gps position wrt to body frame"

# Column E: mirrors column B via formula, consistent with rows above
$ws.Range("E26").Formula = "=B26"
$ws.Range("E27").Formula = "=B27"
$ws.Range("E28").Formula = "=B28"

# --- Update the active selection shown in the sheet view --------------------
[void]$ws.Activate()
$ws.Range("B11").Select() | Out-Null

Write-Host "GPS offset rows added"
